$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue "D2" "42.996.30"
Set-TextValue "E2" "  +0.17%  "
Set-TextValue "D3" "2.540.64"
Set-TextValue "E3" "  -0.78%  "
Set-TextValue "E4" "  -0.15%  "
Set-TextValue "D5" "306.31"
Set-TextValue "E5" "  +1.44%  "
Set-TextValue "D6" "100.90"
Set-TextValue "E6" "  +9.06%  "
Set-TextValue "E7" "  +1.14%  "
Set-TextValue "E8" "  +0.04%  "
Set-TextValue "D9" "0.550"
Set-TextValue "E9" "  +0.78%  "
Set-TextValue "E10" "  +4.04%  "
Set-TextValue "E11" "  +1.20%  "
Set-TextValue "D12" "7.66"
Set-TextValue "E12" "  -0.53%  "
Set-TextValue "E13" "  -0.40%  "
Set-TextValue "D14" "2.929.41"
Set-TextValue "E14" "  -0.89%  "
Set-TextValue "D15" "2.541.10"
Set-TextValue "E15" "  +0.87%  "
Set-TextValue "D16" "15.33"
Set-TextValue "E16" "  +7.94%  "
Set-TextValue "D17" "0.874"
Set-TextValue "E17" "  -0.73%  "
Set-TextValue "D18" "42.980.41"
Set-TextValue "E18" "  +0.01%  "
Set-TextValue "D19" "13.15"
Set-TextValue "E19" "  +3.73%  "
Set-TextValue "D20" "0.0₃0990"
Set-TextValue "E20" "  -0.46%  "
Set-TextValue "D21" "6.51"
Set-TextValue "D22" "71.76"
Set-TextValue "E22" "  +0.25%  "
Set-TextValue "D23" "254.60"
Set-TextValue "E23" "  +0.54%  "
Set-TextValue "D24" "2.94"
Set-TextValue "E24" "  -0.06%  "
Set-TextValue "E25" "  -3.09%  "
Set-TextValue "D26" "27.46"
Set-TextValue "E26" "  -4.32%  "
Set-TextValue "E27" "  +0.27%  "
Set-TextValue "D28" "10.50"
Set-TextValue "E28" "  +2.07%  "
Set-TextValue "D29" "2.34"
Set-TextValue "E29" "  +9.76%  "
Set-TextValue "D30" "39.35"
Set-TextValue "E30" "  +5.90%  "
Set-TextValue "E31" "  +2.19%  "
Set-TextValue "D32" "158.69"
Set-TextValue "E32" "  +3.19%  "
Set-TextValue "B33" "Hedera"
Set-TextValue "C33" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D33" "0.0802"
Set-TextValue "E33" "  +0.33%  "
Set-TextValue "B34" "ARBITRUM"
Set-TextValue "C34" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D34" "2.12"
Set-TextValue "E34" "  -1.23%  "
Set-TextValue "D35" "3.30"
Set-TextValue "E35" "  -2.25%  "
Set-TextValue "E36" "  -3.61%  "
Set-TextValue "D37" "18.57"
Set-TextValue "E37" "  +3.12%  "
Set-TextValue "E38" "  +1.65%  "
Set-TextValue "B39" "Stellar"
Set-TextValue "C39" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D39" "0.120"
Set-TextValue "E39" "  +0.22%  "
Set-TextValue "B40" "EnergySwap"
Set-TextValue "C40" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D40" "24.18"
Set-TextValue "E40" "  +4.93%  "
Set-TextValue "E41" "  +1.57%  "
Set-TextValue "E42" "  +3.65%  "
Set-TextValue "D43" "3.90"
Set-TextValue "E43" "  +0.69%  "
Set-TextValue "E44" "  -1.71%  "
Set-TextValue "E45" "  -0.02%  "
Set-TextValue "D46" "2.053.19"
Set-TextValue "E46" "  -2.09%  "
Set-TextValue "D47" "86.27"
Set-TextValue "E47" "  +1.44%  "
Set-TextValue "E48" "  -2.73%  "
Set-TextValue "D49" "2.786.14"
Set-TextValue "E49" "  -0.89%  "
Set-TextValue "E50" "  +1.42%  "
Set-TextValue "D51" "103.87"
Set-TextValue "E51" "  -2.50%  "
